$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 701, shifting existing rows 701:795 down to 702:796
$ws.Rows.Item(701).Insert()

# Populate the newly inserted row 701 with its data
$ws.Cells.Item(701, 1).Value = 5
$ws.Cells.Item(701, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(701, 3).Value = "Maule"
$ws.Cells.Item(701, 4).Value = 44984
$ws.Cells.Item(701, 5).Value = 7
$ws.Cells.Item(701, 6).Value = "Fruta"
$ws.Cells.Item(701, 7).Value = 100102
$ws.Cells.Item(701, 8).Value = "Cítricos"
$ws.Cells.Item(701, 9).Value = 100102005
$ws.Cells.Item(701, 10).Value = "Naranja"
$ws.Cells.Item(701, 11).Value = "Valencia"
$ws.Cells.Item(701, 12).Value = "Primera"
$ws.Cells.Item(701, 13).Value = 200
$ws.Cells.Item(701, 14).Value = 14000
$ws.Cells.Item(701, 15).Value = 14000
$ws.Cells.Item(701, 16).Value = 14000
$ws.Cells.Item(701, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(701, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(701, 19).Value = 933
$ws.Cells.Item(701, 20).Value = 15
